$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.910.66'
$ws.Range('E2').Value = '  +1.37%  '

$ws.Range('D3').Value = '2.701.35'
$ws.Range('E3').Value = '  +2.49%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '''609.02'
$ws.Range('E5').Value = '  +2.02%  '

$ws.Range('D6').Value = '''158.15'
$ws.Range('E6').Value = '  +1.49%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('E8').Value = '  -0.53%  '

$ws.Range('E9').Value = '  +5.68%  '

$ws.Range('D10').Value = '''6.01'
$ws.Range('E10').Value = '  +3.80%  '

$ws.Range('D11').Value = '''0.403'
$ws.Range('E11').Value = '  +0.58%  '

$ws.Range('E12').Value = '  +0.82%  '

$ws.Range('D13').Value = '''30.46'
$ws.Range('E13').Value = '  +4.31%  '

$ws.Range('D14').Value = '''0.0000203'
$ws.Range('E14').Value = '  +8.97%  '

$ws.Range('D15').Value = '3.185.56'
$ws.Range('E15').Value = '  +2.43%  '

$ws.Range('D16').Value = '65.765.59'
$ws.Range('E16').Value = '  +1.28%  '

$ws.Range('D17').Value = '2.702.28'
$ws.Range('E17').Value = '  +2.46%  '

$ws.Range('D18').Value = '''12.70'
$ws.Range('E18').Value = '  +1.26%  '

$ws.Range('D19').Value = '''4.90'
$ws.Range('E19').Value = '  +2.17%  '

$ws.Range('D20').Value = '''359.42'
$ws.Range('E20').Value = '  +1.99%  '

$ws.Range('D21').Value = '''7.58'
$ws.Range('E21').Value = '  +3.45%  '

$ws.Range('D22').Value = '''1.00'
$ws.Range('E22').Value = '  -0.11%  '

$ws.Range('D23').Value = '''70.81'
$ws.Range('E23').Value = '  +3.77%  '

$ws.Range('D24').Value = '''9.88'
$ws.Range('E24').Value = '  +3.63%  '

$ws.Range('D25').Value = '''0.0000107'
$ws.Range('E25').Value = '  +13.03%  '

$ws.Range('D26').Value = '''1.67'
$ws.Range('E26').Value = '  -1.62%  '

$ws.Range('D27').Value = '''1.69'
$ws.Range('E27').Value = '  +2.95%  '

$ws.Range('E28').Value = '  +4.13%  '

$ws.Range('D29').Value = '''8.41'
$ws.Range('E29').Value = '  +4.11%  '

$ws.Range('E30').Value = '  +4.82%  '

$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  +0.06%  '

$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '''539.64'
$ws.Range('E32').Value = '  +6.20%  '

$ws.Range('D33').Value = '''1.81'
$ws.Range('E33').Value = '  +2.71%  '

$ws.Range('D34').Value = '''6.76'
$ws.Range('E34').Value = '  +6.70%  '

$ws.Range('D35').Value = '''5.43'
$ws.Range('E35').Value = '  -2.84%  '

$ws.Range('D36').Value = '''0.434'
$ws.Range('E36').Value = '  +2.14%  '

$ws.Range('D37').Value = '''20.91'
$ws.Range('E37').Value = '  +3.12%  '

$ws.Range('D38').Value = '''162.89'
$ws.Range('E38').Value = '  -0.26%  '

$ws.Range('D39').Value = '''2.00'
$ws.Range('E39').Value = '  -0.01%  '

$ws.Range('E40').Value = '  +0.05%  '

$ws.Range('D41').Value = '''172.06'
$ws.Range('E41').Value = '  +3.97%  '

$ws.Range('E42').Value = '  +0.04%  '

$ws.Range('D43').Value = '''42.45'
$ws.Range('E43').Value = '  +0.40%  '

$ws.Range('D44').Value = '''4.19'
$ws.Range('E44').Value = '  +2.84%  '

$ws.Range('D45').Value = '''0.0619'
$ws.Range('E45').Value = '  +0.33%  '

$ws.Range('D46').Value = '''23.61'
$ws.Range('E46').Value = '  +2.21%  '

$ws.Range('D47').Value = '''2.29'
$ws.Range('E47').Value = '  +3.84%  '

$ws.Range('D48').Value = '''0.0267'
$ws.Range('E48').Value = '  +4.51%  '

$ws.Range('D49').Value = '''0.655'
$ws.Range('E49').Value = '  +1.37%  '

$ws.Range('D50').Value = '''21.12'
$ws.Range('E50').Value = '  +8.45%  '

$ws.Range('D51').Value = '''0.0992'
$ws.Range('E51').Value = '  +0.93%  '
